$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header label for the hybrid column in the first table (row 2)
$ws.Range("F2").Value = "Hybrid"

# Add header row for the second (hybrid) mini-table starting at row 40,
# mirroring the labels already used in row 2 (B:E)
$ws.Range("B40").Value = $ws.Range("B2").Value()
$ws.Range("C40").Value = $ws.Range("C2").Value()
$ws.Range("D40").Value = $ws.Range("D2").Value()
$ws.Range("E40").Value = $ws.Range("E2").Value()

# Fill in the new threshold / optimization values for rows 41-47
$ws.Range("B41").Value = 1000
$ws.Range("C41").Value = 257
$ws.Range("D41").Value = 373
$ws.Range("E41").Value = 383

$ws.Range("B42").Value = 3000
$ws.Range("C42").Value = 638
$ws.Range("D42").Value = 533
$ws.Range("E42").Value = 736

$ws.Range("B43").Value = 10000
$ws.Range("C43").Value = 2683
$ws.Range("D43").Value = 1162
$ws.Range("E43").Value = 3247

$ws.Range("B44").Value = 30000
$ws.Range("C44").Value = 18726
$ws.Range("D44").Value = 3231
$ws.Range("E44").Value = 23968

$ws.Range("B45").Value = 100000
$ws.Range("C45").Value = 217040
$ws.Range("D45").Value = 15731
$ws.Range("E45").Value = 259537

$ws.Range("B46").Value = 300000

$ws.Range("B47").Value = 1000000

# Update the view: clear the scrolled top-left cell and move the active
# selection down to C53
$ws.Range("C53").Select() | Out-Null
